# Auto-generated edit script: refresh crypto price/volume snapshot
# (rank column A is untouched; coin/link/price/volume columns updated
# to match the latest scrape, including the list reshuffling where
# "Frax" drops off the board and "Aave" is newly appended at the end)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.532.43"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.732.87"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4922"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2638"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06187"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").Value = "1.739.31"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07026"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.564"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6025"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "26.538.76"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007152"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("D21").Value = "1.963.27"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.504"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.616"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.202"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.443"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "107.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.727"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.979"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07985"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.679"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04537"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.616"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.005"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.72%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6271"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9084"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.019"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.18%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.411"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.003"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01493"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.95%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.462"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.08%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3881"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.707"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.93%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1161"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05368"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.63%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.719"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.251"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.01%  "
